$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.764.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'2.103.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'228.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'62.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.72%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").Value = "'0.0842"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "'15.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.44%  "
$ws.Range("D13").Value = "'2.415.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'22.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "'0.809"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "'5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "'2.100.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'38.807.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'71.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "'6.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "'0.0₃0840"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'228.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.15%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'171.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  +5.83%  "
$ws.Range("D29").Value = "'1.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.41%  "
$ws.Range("D30").Value = "'19.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("E31").Value = "  +8.08%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "'6.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.21%  "
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "'3.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'18.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "'102.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").Value = "'1.533.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "'7.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'2.301.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "
